$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- General info text updates ---
$ws.Range("A5").Value = "Issue date: 14/01/2021 16:03:45"
$ws.Range("A6").Value = "Python version: Python 3.8.5"

# --- Row 16 (single remaining data row) updates ---
$ws.Range("A16").Value = 1001
$ws.Range("C16").Value = "Stylistic Features: vof, huf, aof, pnf, anf, agf, frc, mef, acf, fdf"
$ws.Range("F16").Value = "5 folds X 20 iterations CV"
$ws.Range("J16").Value = "90.15V"

# Re-style J16 to match the "V" (significantly larger) red/centered style used
# elsewhere in the workbook (same style previously used on J18).
$ws.Range("F11").Copy()
$ws.Range("J16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Remove the now-obsolete data rows (old rows 17-20) ---
$ws.Rows("17:20").Delete()

# --- Column C width change ---
$ws.Columns.Item(3).ColumnWidth = 67.75

# --- Table style rename ---
$tbl = $ws.ListObjects.Item(1)
$tbl.TableStyle = "TableStyleLight13"
